$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 holds the single data record; update it to the new tax-record values.
# All cells in this sheet store plain text (shared strings), so purely
# numeric-looking values are entered with a leading apostrophe - the normal
# Excel way of forcing text entry - to keep them as text instead of letting
# Excel reinterpret them as numbers.
$ws.Range("B3").Value = "14/07/2022"
$ws.Range("C3").Value = "MIGROS TICARET A.S."
$ws.Range("D3").Value = "'0506"
$ws.Range("E3").Value = "BUYUK MUKELLEFLER "
$ws.Range("F3").Value = "'6220529513"
$ws.Range("G3").Value = "'31.11"
$ws.Range("H3").Value = "'172.83"
$ws.Range("I3").Value = "'203.94"
